$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto data (prices / volume %) as scraped on the new run.
# All target cells are plain text cells (inlineStr in the original workbook),
# so values are assigned as strings to avoid Excel auto-converting them to
# numbers/dates and to preserve the exact padding/format of the percentages.

# Row 2
$ws.Cells.Item(2, 4).Value = "72.484.82"
$ws.Cells.Item(2, 5).Value = "  +2.19%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.643.69"
$ws.Cells.Item(3, 5).Value = "  +0.90%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.04%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "603.82"
$ws.Cells.Item(5, 5).Value = "  -0.16%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "180.29"
$ws.Cells.Item(6, 5).Value = "  -0.18%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  -0.03%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +0.32%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  +7.21%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "2.643.24"
$ws.Cells.Item(10, 5).Value = "  +0.91%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "0.358"
$ws.Cells.Item(12, 5).Value = "  +3.51%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "5.05"
$ws.Cells.Item(13, 5).Value = "  +0.41%  "

# Row 14
$ws.Cells.Item(14, 5).Value = "  +4.85%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "3.127.78"
$ws.Cells.Item(15, 5).Value = "  +1.27%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "72.327.86"
$ws.Cells.Item(16, 5).Value = "  +1.98%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "26.57"
$ws.Cells.Item(17, 5).Value = "  -0.39%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "2.648.58"
$ws.Cells.Item(18, 5).Value = "  +1.46%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "  +4.15%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "379.84"
$ws.Cells.Item(20, 5).Value = "  -0.19%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "7.93"
$ws.Cells.Item(21, 5).Value = "  +0.21%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "4.18"
$ws.Cells.Item(22, 5).Value = "  +0.46%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  +11.20%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "72.99"
$ws.Cells.Item(24, 5).Value = "  +1.51%  "

# Row 25
$ws.Cells.Item(25, 2).Value = "Dai"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(25, 4).Value = "1.00"
$ws.Cells.Item(25, 5).Value = "  -0.17%  "

# Row 26
$ws.Cells.Item(26, 2).Value = "NEARProtocol"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(26, 4).Value = "4.39"
$ws.Cells.Item(26, 5).Value = "  -0.98%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  +3.84%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "2.780.00"
$ws.Cells.Item(28, 5).Value = "  +1.06%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  +0.05%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "0.0₃0961"
$ws.Cells.Item(30, 5).Value = "  +1.62%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "524.89"
$ws.Cells.Item(31, 5).Value = "  -0.06%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "8.15"
$ws.Cells.Item(32, 5).Value = "  +1.53%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "1.31"
$ws.Cells.Item(33, 5).Value = "  -0.70%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "1.82"
$ws.Cells.Item(34, 5).Value = "  -0.30%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "0.999"

# Row 36
$ws.Cells.Item(36, 4).Value = "164.81"
$ws.Cells.Item(36, 5).Value = "  +0.04%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "19.32"
$ws.Cells.Item(37, 5).Value = "  +1.08%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  -5.54%  "

# Row 39
$ws.Cells.Item(39, 4).Value = "19.10"
$ws.Cells.Item(39, 5).Value = "  +0.88%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  +2.04%  "

# Row 41
$ws.Cells.Item(41, 4).Value = "1.87"
$ws.Cells.Item(41, 5).Value = "  -0.31%  "

# Row 42
$ws.Cells.Item(42, 2).Value = "RenderToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Cells.Item(42, 4).Value = "5.10"
$ws.Cells.Item(42, 5).Value = "  +1.46%  "

# Row 43
$ws.Cells.Item(43, 2).Value = "dogwifhat"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(43, 4).Value = "2.64"
$ws.Cells.Item(43, 5).Value = "  +1.97%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  -0.01%  "

# Row 45
$ws.Cells.Item(45, 5).Value = "  +0.68%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "39.39"
$ws.Cells.Item(46, 5).Value = "  -1.80%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "151.44"
$ws.Cells.Item(47, 5).Value = "  -1.37%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  +2.10%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "  +2.34%  "

# Row 50
$ws.Cells.Item(50, 5).Value = "  +2.71%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "0.0₆0260"
$ws.Cells.Item(51, 5).Value = "  -3.82%  "
